$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row: "Locations" -> "probeCoords", "ProbeLocation" -> "stimOverProbe"
$ws.Range("C1").Value = "probeCoords"
$ws.Range("D1").Value = "stimOverProbe"

# Update existing data rows (2 and 3) to use the [0.35, 0] coordinate for both,
# and change the stimOverProbe label to reflect which face is on top (opposite of before)
$ws.Range("C2").Value = "[0.35, 0]"
$ws.Range("D2").Value = "frown"

$ws.Range("C3").Value = "[0.35, 0]"
$ws.Range("D3").Value = "smile"

# Add two new rows (4 and 5) with the [-0.35, 0] coordinate
$ws.Range("A4").Value = "smile.jpg"
$ws.Range("B4").Value = "frown.jpg"
$ws.Range("C4").Value = "[-0.35, 0]"
$ws.Range("D4").Value = "frown"

$ws.Range("A5").Value = "frown.jpg"
$ws.Range("B5").Value = "smile.jpg"
$ws.Range("C5").Value = "[-0.35, 0]"
$ws.Range("D5").Value = "smile"

# Update the selected cell in the sheet view
$ws.Range("F7").Select()
